$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = @('ECs', 'Vegfc', 'Flt4', 'ECs', 3, 1, 2.404594333333333, 7.213783, 0.3565065439253589, 0.3565065439253589, 3, 1, 16.203114, 48.609342, 0.9864083027928519, 0.9864083027928516, 38.961916106754, 350.657244960786, 0.3516610149279585, 0.3516610149279584)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(2, $c).Value = $row[$c-1] }

$row = @('ECs', 'Vegfc', 'Flt4', 'FAPs', 3, 1, 2.404594333333333, 7.213783, 0.3565065439253589, 0.3565065439253589, 3, 1, 0.163766, 0.491298, 0.009969697313440749, 0.009969697313440747, 0.3937907955926667, 3.544117160334, 0.003554262333196696, 0.003554262333196696)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(3, $c).Value = $row[$c-1] }

$row = @('ECs', 'Vegfc', 'Flt4', 'MuSCs', 3, 1, 2.404594333333333, 7.213783, 0.3565065439253589, 0.3565065439253589, 2, 0.6666666666666666, 0.03620466666666667, 0.108614, 0.002204056812773619, 0.002204056812773619, 0.08705753630688889, 0.783517826762, 0.0007857606769370647, 0.0007857606769370645)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(4, $c).Value = $row[$c-1] }

$row = @('ECs', 'Vegfc', 'Flt4', 'Neutrophils', 3, 1, 2.404594333333333, 7.213783, 0.3565065439253589, 0.3565065439253589, 1, 0.3333333333333333, 0.02329166666666667, 0.06987500000000001, 0.001417943080933919, 0.001417943080933918, 0.05600700968055556, 0.504063087125, 0.0005055059872666268, 0.0005055059872666266)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(5, $c).Value = $row[$c-1] }

$row = @('FAPs', 'Vegfc', 'Flt4', 'ECs', 3, 1, 3.178631333333334, 9.535894000000001, 0.471265716362492, 0.4712657163624919, 3, 1, 16.203114, 48.609342, 0.9864083027928519, 0.9864083027928516, 51.503725857972, 463.533532721748, 0.4648604154415832, 0.4648604154415831)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(6, $c).Value = $row[$c-1] }

$row = @('FAPs', 'Vegfc', 'Flt4', 'FAPs', 3, 1, 3.178631333333334, 9.535894000000001, 0.471265716362492, 0.4712657163624919, 3, 1, 0.163766, 0.491298, 0.009969697313440749, 0.009969697313440747, 0.5205517389346667, 4.684965650412001, 0.004698376546335866, 0.004698376546335865)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(7, $c).Value = $row[$c-1] }

$row = @('FAPs', 'Vegfc', 'Flt4', 'MuSCs', 3, 1, 3.178631333333334, 9.535894000000001, 0.471265716362492, 0.4712657163624919, 2, 0.6666666666666666, 0.03620466666666667, 0.108614, 0.002204056812773619, 0.002204056812773619, 0.1150812878795556, 1.035731590916, 0.00103869641277539, 0.00103869641277539)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(8, $c).Value = $row[$c-1] }

$row = @('FAPs', 'Vegfc', 'Flt4', 'Neutrophils', 3, 1, 3.178631333333334, 9.535894000000001, 0.471265716362492, 0.4712657163624919, 1, 0.3333333333333333, 0.02329166666666667, 0.06987500000000001, 0.001417943080933919, 0.001417943080933918, 0.07403562147222223, 0.6663205932500001, 0.0006682279617975621, 0.0006682279617975618)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(9, $c).Value = $row[$c-1] }

$row = @('MuSCs', 'Vegfc', 'Flt4', 'ECs', 3, 1, 1.125649, 3.376947, 0.166889370527102, 0.1668893705271019, 3, 1, 16.203114, 48.609342, 0.9864083027928519, 0.9864083027928516, 18.239019070986, 164.151171638874, 0.164621060735806, 0.164621060735806)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(10, $c).Value = $row[$c-1] }

$row = @('MuSCs', 'Vegfc', 'Flt4', 'FAPs', 3, 1, 1.125649, 3.376947, 0.166889370527102, 0.1668893705271019, 3, 1, 0.163766, 0.491298, 0.009969697313440749, 0.009969697313440747, 0.184343034134, 1.659087307206, 0.001663836508985866, 0.001663836508985866)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(11, $c).Value = $row[$c-1] }

$row = @('MuSCs', 'Vegfc', 'Flt4', 'MuSCs', 3, 1, 1.125649, 3.376947, 0.166889370527102, 0.1668893705271019, 2, 0.6666666666666666, 0.03620466666666667, 0.108614, 0.002204056812773619, 0.002204056812773619, 0.04075374682866668, 0.366783721458, 0.0003678336540897599, 0.0003678336540897598)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(12, $c).Value = $row[$c-1] }

$row = @('MuSCs', 'Vegfc', 'Flt4', 'Neutrophils', 3, 1, 1.125649, 3.376947, 0.166889370527102, 0.1668893705271019, 1, 0.3333333333333333, 0.02329166666666667, 0.06987500000000001, 0.001417943080933919, 0.001417943080933918, 0.02621824129166667, 0.2359641716250001, 0.0002366396282203213, 0.0002366396282203212)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(13, $c).Value = $row[$c-1] }

$row = @('Neutrophils', 'Vegfc', 'Flt4', 'ECs', 1, 0.3333333333333333, 0.03600666666666667, 0.10802, 0.005338369185047189, 0.005338369185047188, 3, 1, 16.203114, 48.609342, 0.9864083027928519, 0.9864083027928516, 0.5834201247599999, 5.25078112284, 0.005265811687504058, 0.005265811687504056)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(14, $c).Value = $row[$c-1] }

$row = @('Neutrophils', 'Vegfc', 'Flt4', 'FAPs', 1, 0.3333333333333333, 0.03600666666666667, 0.10802, 0.005338369185047189, 0.005338369185047188, 3, 1, 0.163766, 0.491298, 0.009969697313440749, 0.009969697313440747, 0.005896667773333333, 0.05307000996, 0.00005322192492231984, 0.00005322192492231982)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(15, $c).Value = $row[$c-1] }

$row = @('Neutrophils', 'Vegfc', 'Flt4', 'MuSCs', 1, 0.3333333333333333, 0.03600666666666667, 0.10802, 0.005338369185047189, 0.005338369185047188, 2, 0.6666666666666666, 0.03620466666666667, 0.108614, 0.002204056812773619, 0.002204056812773619, 0.001303609364444445, 0.01173248428, 0.00001176606897140401, 0.00001176606897140401)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(16, $c).Value = $row[$c-1] }

$row = @('Neutrophils', 'Vegfc', 'Flt4', 'Neutrophils', 1, 0.3333333333333333, 0.03600666666666667, 0.10802, 0.005338369185047189, 0.005338369185047188, 1, 0.3333333333333333, 0.02329166666666667, 0.06987500000000001, 0.001417943080933919, 0.001417943080933918, 0.0008386552777777779, 0.007547897500000001, 0.000007569503649408505, 0.000007569503649408501)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(17, $c).Value = $row[$c-1] }
